# FeedBack Survey Form F9.xlsx edit script
# Renames sheet F-SW-CR-07 -> F-SW-SD-09, updates the print area defined name
# accordingly, updates the active cell selection / scroll position, updates
# the footer revision text, and nudges the workbook window position/size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the main worksheet (Print_Area defined name follows automatically) ---
$ws.Name = "F-SW-SD-09"

# Re-assert the print area so the Print_Area defined name points at the
# renamed sheet with the expected absolute reference.
$ws.PageSetup.PrintArea = "`$A`$1:`$F`$29"

# --- Update the selection / scroll position on the sheet ---
$ws.Range("C18:F18").Select()

# --- Update the footer revision text ---
$ws.PageSetup.RightFooter = "&`"Arial,Regular`"&16Rev:0(01/10/2025)"

# --- Nudge the workbook window geometry ---
$win = $excel.ActiveWindow
$win.Left = 20370
$win.Top = -2595
$win.Width = 29040
$win.Height = 15840
